# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# New K values for rows 2..76 (column G), pulled/recalculated from the
# updated source data.
$kValues = @(2,0,0,1,0,2,1,1,0,0,2,1,0,0,1,1,0,1,1,0,1,0,1,2,0,1,2,2,2,1,1,1,0,1,2,1,0,0,1,0,1,0,3,0,0,1,2,1,2,0,2,1,0,1,3,0,1,1,0,1,0,0,0,2,1,2,1,2,2,2,2,2,1,1,1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
